# Auto-generated edit script: updates crypto price/volume data in Sheet1
# Applies the inline-string cell updates described by the commit diff,
# forcing text storage so numeric-looking strings (e.g. "0.999", "503.34")
# are not silently re-typed as numbers by Excel's input parser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Addr, $Val)
    $cell = $ws.Range($Addr)
    # Force text interpretation so strings like "0.999" or "503.34" aren't
    # auto-converted to numbers, then strip the temporary format again so the
    # cell's style stays exactly as it was (no stray numFmt/quote-prefix left).
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.ClearFormats()
}

Set-TextValue "D2" '56.498.75'
Set-TextValue "E2" '  -2.14%  '
Set-TextValue "D3" '2.381.56'
Set-TextValue "E3" '  -1.13%  '
Set-TextValue "E4" '  +0.04%  '
Set-TextValue "D5" '503.34'
Set-TextValue "E5" '  -0.62%  '
Set-TextValue "D6" '130.68'
Set-TextValue "E6" '  -1.72%  '
Set-TextValue "E7" '  +0.25%  '
Set-TextValue "E8" '  -2.06%  '
Set-TextValue "D9" '2.389.12'
Set-TextValue "E9" '  -2.32%  '
Set-TextValue "E10" '  +1.03%  '
Set-TextValue "E11" '  +0.29%  '
Set-TextValue "E12" '  +1.52%  '
Set-TextValue "D13" '4.74'
Set-TextValue "E13" '  +3.23%  '
Set-TextValue "D14" '2.802.23'
Set-TextValue "E14" '  -1.29%  '
Set-TextValue "D15" '56.444.69'
Set-TextValue "E15" '  -1.41%  '
Set-TextValue "E16" '  -1.07%  '
Set-TextValue "E17" '  -0.36%  '
Set-TextValue "D18" '2.337.05'
Set-TextValue "E18" '  -4.74%  '
Set-TextValue "E19" '  -2.42%  '
Set-TextValue "E20" '  -1.61%  '
Set-TextValue "D21" '307.83'
Set-TextValue "E21" '  -2.01%  '
Set-TextValue "D22" '6.28'
Set-TextValue "E22" '  -1.82%  '
Set-TextValue "D23" '0.999'
Set-TextValue "D24" '65.60'
Set-TextValue "E24" '  +0.23%  '
Set-TextValue "E25" '  +0.19%  '
Set-TextValue "D26" '0.368'
Set-TextValue "E26" '  -3.51%  '
Set-TextValue "E27" '  -3.09%  '
Set-TextValue "D28" '7.32'
Set-TextValue "E28" '  -3.33%  '
Set-TextValue "D29" '172.05'
Set-TextValue "E29" '  -1.13%  '
Set-TextValue "E30" '  -1.72%  '
Set-TextValue "E31" '  -2.29%  '
Set-TextValue "B33" 'Fetch.AI'
Set-TextValue "C33" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D33" '1.09'
Set-TextValue "E33" '  -3.12%  '
Set-TextValue "B34" 'Aptos'
Set-TextValue "C34" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D34" '5.77'
Set-TextValue "E34" '  -6.22%  '
Set-TextValue "D35" '0.996'
Set-TextValue "E35" '  +0.16%  '
Set-TextValue "D36" '17.62'
Set-TextValue "E36" '  -1.70%  '
Set-TextValue "E37" '  -4.88%  '
Set-TextValue "E38" '  -1.27%  '
Set-TextValue "E39" '  -1.13%  '
Set-TextValue "D40" '0.798'
Set-TextValue "E40" '  -1.87%  '
Set-TextValue "E41" '  -4.21%  '
Set-TextValue "D42" '130.98'
Set-TextValue "E42" '  -2.68%  '
Set-TextValue "E43" '  -0.32%  '
Set-TextValue "D44" '4.98'
Set-TextValue "E44" '  +0.01%  '
Set-TextValue "E45" '  -0.84%  '
Set-TextValue "D46" '0.0908'
Set-TextValue "E46" '  -0.78%  '
Set-TextValue "E47" '  -5.54%  '
Set-TextValue "D48" '0.0483'
Set-TextValue "E48" '  -1.87%  '
Set-TextValue "E49" '  -1.75%  '
Set-TextValue "D50" '17.17'
Set-TextValue "E50" '  +0.47%  '
Set-TextValue "E51" '  -2.18%  '
